# edit.ps1
# Applies the diff: update Resumen.C2, Solucion!A2:B41 (reorder of Pedido/Salida assignments),
# and Metricas.B2:B3 (zone Wmax/time metrics) to reflect the evolutionary-method's new best solution.

$wb = $excel.ActiveWorkbook
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsSol = $wb.Worksheets.Item("Solucion")
$wsMetricas = $wb.Worksheets.Item("Metricas")

# --- Resumen: update Maximo (Wmax) for Zona Z1 ---
$wsResumen.Range("C2").Value = 648.2560547086171

# Update Solucion sheet Pedido (A) and Salida (B) columns per the diff
$wsSol.Range("A2").Value = "Pedido_3"
$wsSol.Range("A3").Value = "Pedido_10"
$wsSol.Range("A4").Value = "Pedido_15"
$wsSol.Range("A5").Value = "Pedido_13"
$wsSol.Range("A6").Value = "Pedido_36"
$wsSol.Range("A7").Value = "Pedido_28"
$wsSol.Range("A8").Value = "Pedido_39"
$wsSol.Range("A9").Value = "Pedido_9"
$wsSol.Range("A10").Value = "Pedido_12"
$wsSol.Range("A11").Value = "Pedido_2"
$wsSol.Range("A12").Value = "Pedido_22"
$wsSol.Range("A13").Value = "Pedido_35"
$wsSol.Range("A14").Value = "Pedido_24"
$wsSol.Range("A15").Value = "Pedido_34"
$wsSol.Range("A16").Value = "Pedido_20"
$wsSol.Range("A17").Value = "Pedido_16"
$wsSol.Range("A18").Value = "Pedido_29"
$wsSol.Range("A19").Value = "Pedido_31"
$wsSol.Range("A20").Value = "Pedido_33"
$wsSol.Range("A21").Value = "Pedido_40"
$wsSol.Range("A22").Value = "Pedido_18"
$wsSol.Range("A23").Value = "Pedido_5"
$wsSol.Range("A24").Value = "Pedido_14"
$wsSol.Range("A25").Value = "Pedido_4"
$wsSol.Range("A26").Value = "Pedido_27"
$wsSol.Range("A27").Value = "Pedido_11"
$wsSol.Range("A28").Value = "Pedido_37"
$wsSol.Range("A29").Value = "Pedido_30"
$wsSol.Range("A30").Value = "Pedido_1"
$wsSol.Range("A31").Value = "Pedido_21"
$wsSol.Range("A32").Value = "Pedido_25"
$wsSol.Range("A33").Value = "Pedido_38"
$wsSol.Range("A34").Value = "Pedido_17"
$wsSol.Range("A35").Value = "Pedido_6"
$wsSol.Range("A36").Value = "Pedido_23"
$wsSol.Range("A37").Value = "Pedido_19"
$wsSol.Range("A38").Value = "Pedido_7"
$wsSol.Range("A39").Value = "Pedido_32"
$wsSol.Range("A40").Value = "Pedido_8"
$wsSol.Range("A41").Value = "Pedido_26"
$wsSol.Range("B4").Value = "S029"
$wsSol.Range("B5").Value = "S005"
$wsSol.Range("B6").Value = "S002"
$wsSol.Range("B7").Value = "S026"
$wsSol.Range("B8").Value = "S006"
$wsSol.Range("B9").Value = "S030"
$wsSol.Range("B10").Value = "S003"
$wsSol.Range("B11").Value = "S027"
$wsSol.Range("B12").Value = "S007"
$wsSol.Range("B13").Value = "S031"
$wsSol.Range("B14").Value = "S004"
$wsSol.Range("B15").Value = "S008"
$wsSol.Range("B19").Value = "S010"
$wsSol.Range("B20").Value = "S032"
$wsSol.Range("B23").Value = "S011"
$wsSol.Range("B24").Value = "S037"
$wsSol.Range("B25").Value = "S015"
$wsSol.Range("B26").Value = "S034"
$wsSol.Range("B27").Value = "S012"
$wsSol.Range("B28").Value = "S038"
$wsSol.Range("B29").Value = "S016"
$wsSol.Range("B30").Value = "S035"

# --- Metricas: update per-zone Tiempo (Wmax) values ---
$wsMetricas.Range("B2").Value = 648.2560547086171
$wsMetricas.Range("B3").Value = 480.6931560168667
